$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.406955
$ws.Range("H2").Value = 19.220865
$ws.Range("I2").Value = 0.2800966009992834
$ws.Range("J2").Value = 0.3266544289500553
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 12.55301866666667
$ws.Range("N2").Value = 37.65905600000001
$ws.Range("O2").Value = 0.2692278546031339
$ws.Range("P2").Value = 0.3068924341132449
$ws.Range("Q2").Value = 80.42662571149334
$ws.Range("R2").Value = 723.8396314034401
$ws.Range("S2").Value = 0.07540980696866709
$ws.Range("T2").Value = 0.1002477728143545
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.406955
$ws.Range("H3").Value = 19.220865
$ws.Range("I3").Value = 0.2800966009992834
$ws.Range("J3").Value = 0.3266544289500553
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.354491666666666
$ws.Range("N3").Value = 28.063475
$ws.Range("O3").Value = 0.200628214551068
$ws.Range("P3").Value = 0.2286958056629511
$ws.Range("Q3").Value = 59.93380715620833
$ws.Range("R3").Value = 539.404264405875
$ws.Range("S3").Value = 0.05619528096030911
$ws.Range("T3").Value = 0.07470449780210413
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.406955
$ws.Range("H4").Value = 19.220865
$ws.Range("I4").Value = 0.2800966009992834
$ws.Range("J4").Value = 0.3266544289500553
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.616516333333333
$ws.Range("N4").Value = 13.849549
$ws.Range("O4").Value = 0.09901162590190733
$ws.Range("P4").Value = 0.1128632062359889
$ws.Range("Q4").Value = 29.57781240443166
$ws.Range("R4").Value = 266.200311639885
$ws.Range("S4").Value = 0.02773281987453685
$ws.Range("T4").Value = 0.03686726618248928
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 6.406955
$ws.Range("H5").Value = 19.220865
$ws.Range("I5").Value = 0.2800966009992834
$ws.Range("J5").Value = 0.3266544289500553
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.934898
$ws.Range("N5").Value = 8.804694
$ws.Range("O5").Value = 0.06294552035656671
$ws.Range("P5").Value = 0.07175150575421439
$ws.Range("Q5").Value = 18.80375941559
$ws.Range("R5").Value = 169.23383474031
$ws.Range("S5").Value = 0.01763082630000554
$ws.Range("T5").Value = 0.02343794713844951
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 6.406955
$ws.Range("H6").Value = 19.220865
$ws.Range("I6").Value = 0.2800966009992834
$ws.Range("J6").Value = 0.3266544289500553
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 17.167078
$ws.Range("N6").Value = 34.334156
$ws.Range("O6").Value = 0.3681867845873241
$ws.Range("P6").Value = 0.2797970482336007
$ws.Range("Q6").Value = 109.98869622749
$ws.Range("R6").Value = 659.93217736494
$ws.Range("S6").Value = 0.1031278668957648
$ws.Range("T6").Value = 0.09139694501265794
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 6.686451000000001
$ws.Range("H7").Value = 20.059353
$ws.Range("I7").Value = 0.2923154911886005
$ws.Range("J7").Value = 0.3409043505233807
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 12.55301866666667
$ws.Range("N7").Value = 37.65905600000001
$ws.Range("O7").Value = 0.2692278546031339
$ws.Range("P7").Value = 0.3068924341132449
$ws.Range("Q7").Value = 83.93514421675202
$ws.Range("R7").Value = 755.4162979507682
$ws.Range("S7").Value = 0.07869947255996822
$ws.Range("T7").Value = 0.1046209659319152
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 6.686451000000001
$ws.Range("H8").Value = 20.059353
$ws.Range("I8").Value = 0.2923154911886005
$ws.Range("J8").Value = 0.3409043505233807
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 9.354491666666666
$ws.Range("N8").Value = 28.063475
$ws.Range("O8").Value = 0.200628214551068
$ws.Range("P8").Value = 0.2286958056629511
$ws.Range("Q8").Value = 62.54835015907501
$ws.Range("R8").Value = 562.9351514316751
$ws.Range("S8").Value = 0.05864673508278736
$ws.Range("T8").Value = 0.07796339509694965
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 6.686451000000001
$ws.Range("H9").Value = 20.059353
$ws.Range("I9").Value = 0.2923154911886005
$ws.Range("J9").Value = 0.3409043505233807
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.616516333333333
$ws.Range("N9").Value = 13.849549
$ws.Range("O9").Value = 0.09901162590190733
$ws.Range("P9").Value = 0.1128632062359889
$ws.Range("Q9").Value = 30.868110253533
$ws.Range("R9").Value = 277.812992281797
$ws.Range("S9").Value = 0.028942632058898
$ws.Range("T9").Value = 0.03847555801986618
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 6.686451000000001
$ws.Range("H10").Value = 20.059353
$ws.Range("I10").Value = 0.2923154911886005
$ws.Range("J10").Value = 0.3409043505233807
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.934898
$ws.Range("N10").Value = 8.804694
$ws.Range("O10").Value = 0.06294552035656671
$ws.Range("P10").Value = 0.07175150575421439
$ws.Range("Q10").Value = 19.624051666998
$ws.Range("R10").Value = 176.616465002982
$ws.Range("S10").Value = 0.01839995070115185
$ws.Range("T10").Value = 0.02446040046821507
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 6.686451000000001
$ws.Range("H11").Value = 20.059353
$ws.Range("I11").Value = 0.2923154911886005
$ws.Range("J11").Value = 0.3409043505233807
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 17.167078
$ws.Range("N11").Value = 34.334156
$ws.Range("O11").Value = 0.3681867845873241
$ws.Range("P11").Value = 0.2797970482336007
$ws.Range("Q11").Value = 114.786825860178
$ws.Range("R11").Value = 688.7209551610681
$ws.Range("S11").Value = 0.1076267007857951
$ws.Range("T11").Value = 0.09538403100643468
# Row 12
$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 9.780684500000001
$ws.Range("H12").Value = 19.561369
$ws.Range("I12").Value = 0.4275879078121161
$ws.Range("J12").Value = 0.332441220526564
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 12.55301866666667
$ws.Range("N12").Value = 37.65905600000001
$ws.Range("O12").Value = 0.2692278546031339
$ws.Range("P12").Value = 0.3068924341132449
$ws.Range("Q12").Value = 122.7771151012774
$ws.Range("R12").Value = 736.6626906076642
$ws.Range("S12").Value = 0.1151185750744986
$ws.Range("T12").Value = 0.1020236953669753
# Row 13
$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 9.780684500000001
$ws.Range("H13").Value = 19.561369
$ws.Range("I13").Value = 0.4275879078121161
$ws.Range("J13").Value = 0.332441220526564
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 9.354491666666666
$ws.Range("N13").Value = 28.063475
$ws.Range("O13").Value = 0.200628214551068
$ws.Range("P13").Value = 0.2286958056629511
$ws.Range("Q13").Value = 91.49333164954584
$ws.Range("R13").Value = 548.9599898972751
$ws.Range("S13").Value = 0.08578619850797149
$ws.Range("T13").Value = 0.07602791276389737
# Row 14
$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 9.780684500000001
$ws.Range("H14").Value = 19.561369
$ws.Range("I14").Value = 0.4275879078121161
$ws.Range("J14").Value = 0.332441220526564
$ws.Range("K14").Value = 2
$ws.Range("M14").Value = 4.616516333333333
$ws.Range("N14").Value = 13.849549
$ws.Range("O14").Value = 0.09901162590190733
$ws.Range("P14").Value = 0.1128632062359889
$ws.Range("Q14").Value = 45.15268974543017
$ws.Range("R14").Value = 270.9161384725811
$ws.Range("S14").Value = 0.04233617396847248
$ws.Range("T14").Value = 0.03752038203363347
# Row 15
$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 9.780684500000001
$ws.Range("H15").Value = 19.561369
$ws.Range("I15").Value = 0.4275879078121161
$ws.Range("J15").Value = 0.332441220526564
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 2.934898
$ws.Range("N15").Value = 8.804694
$ws.Range("O15").Value = 0.06294552035656671
$ws.Range("P15").Value = 0.07175150575421439
$ws.Range("Q15").Value = 28.705311377681
$ws.Range("R15").Value = 172.231868266086
$ws.Range("S15").Value = 0.02691474335540932
$ws.Range("T15").Value = 0.02385315814754981
# Row 16
$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 9.780684500000001
$ws.Range("H16").Value = 19.561369
$ws.Range("I16").Value = 0.4275879078121161
$ws.Range("J16").Value = 0.332441220526564
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 17.167078
$ws.Range("N16").Value = 34.334156
$ws.Range("O16").Value = 0.3681867845873241
$ws.Range("P16").Value = 0.2797970482336007
$ws.Range("Q16").Value = 167.905773704891
$ws.Range("R16").Value = 671.6230948195641
$ws.Range("S16").Value = 0.1574322169057642
$ws.Range("T16").Value = 0.09301607221450812

Write-Host "Applied 210 cell updates"